$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (matches upstream coinranking.com scrape refresh).
# D-column (Price) cells are forced to Text via NumberFormat "@" so numeric-looking
# strings (e.g. "61.60", "87.40") keep their exact trailing-zero text form instead of
# being auto-coerced to a Number by the COM Value setter; Style is reset to "Normal"
# afterwards so no stray per-cell formatting is introduced (matches original: no `s` attr).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.911.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.549.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.36%  '

$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("E7").Value = '  -0.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.12'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.94%  '

$ws.Range("E9").Value = '  -0.66%  '

$ws.Range("E10").Value = '  +0.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0855'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.770.85'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.547.38'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.55%  '

$ws.Range("E14").Value = '  +0.63%  '

$ws.Range("E15").Value = '  +0.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.902.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.18'
$ws.Range("D18").Style = "Normal"

$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.30%  '

$ws.Range("E24").Value = '  -0.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.31%  '

$ws.Range("E28").Value = '  +0.65%  '

$ws.Range("E29").Value = '  -0.36%  '

$ws.Range("E30").Value = '  +1.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.05%  '

$ws.Range("E32").Value = '  -0.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.422.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.85%  '

$ws.Range("E34").Value = '  +4.03%  '

$ws.Range("E35").Value = '  +2.37%  '

$ws.Range("E36").Value = '  -0.50%  '

$ws.Range("E37").Value = '  +0.11%  '

$ws.Range("E38").Value = '  +0.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.524'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.14%  '

$ws.Range("E40").Value = '  -0.17%  '

$ws.Range("E41").Value = '  -0.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.71'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.44%  '

$ws.Range("E43").Value = '  +3.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.684.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0960'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.54%  '

